# Updates the cryptocurrency price/volume snapshot in the active worksheet
# to reflect the latest values pulled by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.498.67'
$ws.Range("E2").Value = '  +1.87%  '
$ws.Range("D3").Value = '1.786.29'
$ws.Range("E3").Value = '  -0.03%  '
$ws.Range("E4").Value = '  -0.22%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '222.01'
$c.ClearFormats()
$ws.Range("E5").Value = '  -1.70%  '
$ws.Range("E6").Value = '  -0.75%  '
$ws.Range("E7").Value = '  -0.26%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '32.46'
$c.ClearFormats()
$ws.Range("E8").Value = '  +7.69%  '
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("E10").Value = '  +2.64%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0934'
$c.ClearFormats()
$ws.Range("E11").Value = '  +0.92%  '
$ws.Range("D12").Value = '2.040.92'
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '11.02'
$c.ClearFormats()
$ws.Range("E13").Value = '  +5.43%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.782.36'
$ws.Range("E14").Value = '  -0.28%  '
$ws.Range("D15").Value = '34.490.90'
$ws.Range("E15").Value = '  +1.78%  '
$ws.Range("E16").Value = '  +0.74%  '
$ws.Range("E17").Value = '  +2.18%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '68.68'
$c.ClearFormats()
$ws.Range("E18").Value = '  -0.49%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '253.91'
$c.ClearFormats()
$ws.Range("E19").Value = '  +0.84%  '
$ws.Range("D20").Value = '0.0₃0779'
$ws.Range("E20").Value = '  +5.36%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range("E21").Value = '  -0.21%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '10.49'
$c.ClearFormats()
$ws.Range("E22").Value = '  +1.64%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '4.17'
$c.ClearFormats()
$ws.Range("E23").Value = '  -1.20%  '
$ws.Range("E24").Value = '  +0.16%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '160.34'
$c.ClearFormats()
$ws.Range("E25").Value = '  +1.37%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '16.37'
$c.ClearFormats()
$ws.Range("E26").Value = '  -0.49%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '7.11'
$c.ClearFormats()
$ws.Range("E27").Value = '  +1.74%  '
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("E29").Value = '  -0.29%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.0519'
$c.ClearFormats()
$ws.Range("E30").Value = '  +0.82%  '
$ws.Range("E31").Value = '  -1.95%  '
$ws.Range("E32").Value = '  -0.45%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '3.56'
$c.ClearFormats()
$ws.Range("E33").Value = '  -0.83%  '
$ws.Range("E34").Value = '  +1.06%  '
$ws.Range("D35").Value = '1.430.86'
$ws.Range("E35").Value = '  -4.68%  '
$ws.Range("E36").Value = '  +0.98%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.0191'
$c.ClearFormats()
$ws.Range("E37").Value = '  +3.40%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '1.05'
$c.ClearFormats()
$ws.Range("E38").Value = '  -1.29%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '85.41'
$c.ClearFormats()
$ws.Range("E39").Value = '  +2.48%  '
$ws.Range("E40").Value = '  +3.15%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.923'
$c.ClearFormats()
$ws.Range("E41").Value = '  +2.55%  '
$ws.Range("E42").Value = '  -0.92%  '
$ws.Range("E43").Value = '  +1.84%  '
$ws.Range("E44").Value = '  +4.55%  '
$ws.Range("E45").Value = '  -1.25%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.0492'
$c.ClearFormats()
$ws.Range("E46").Value = '  -5.05%  '
$ws.Range("D47").Value = '1.941.87'
$ws.Range("E47").Value = '  +0.29%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '12.08'
$c.ClearFormats()
$ws.Range("E48").Value = '  +2.27%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '104.32'
$c.ClearFormats()
$ws.Range("E49").Value = '  +6.49%  '
$ws.Range("E50").Value = '  -0.35%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '50.05'
$c.ClearFormats()
$ws.Range("E51").Value = '  -2.34%  '
